# Updates cryptos list values (price + 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.291.69"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "'3.390.73"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'565.34"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'155.62"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'3.391.70"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -8.55%  "
$ws.Range("D10").Value = "'7.22"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  -4.17%  "
$ws.Range("D12").Value = "'0.420"
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("D13").Value = "'3.986.64"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'26.73"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("E16").Value = "  -9.81%  "
$ws.Range("D17").Value = "'63.406.26"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "'3.397.53"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'6.04"
$ws.Range("E19").Value = "  -5.30%  "
$ws.Range("D20").Value = "'13.43"
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("D21").Value = "'381.14"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "'7.69"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'70.91"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").Value = "'0.511"
$ws.Range("E25").Value = "  -7.76%  "
$ws.Range("D26").Value = "'0.0000112"
$ws.Range("E26").Value = "  -5.26%  "
$ws.Range("D27").Value = "'9.62"
$ws.Range("E27").Value = "  -6.26%  "
$ws.Range("D28").Value = "'0.177"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "'5.98"
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("D31").Value = "'1.38"
$ws.Range("E31").Value = "  -8.08%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "'22.71"
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("D35").Value = "'6.86"
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").Value = "'1.49"
$ws.Range("E36").Value = "  -7.78%  "
$ws.Range("D37").Value = "'160.16"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'0.837"
$ws.Range("E38").Value = "  +8.70%  "
$ws.Range("D39").Value = "'1.80"
$ws.Range("E39").Value = "  -5.36%  "
$ws.Range("D40").Value = "'2.807.47"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "'25.79"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Value = "'42.88"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "'0.0710"
$ws.Range("E43").Value = "  -6.81%  "
$ws.Range("D44").Value = "'6.34"
$ws.Range("E44").Value = "  -9.19%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'4.31"
$ws.Range("E45").Value = "  -6.86%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'25.32"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").Value = "'0.0301"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'323.17"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.28"
$ws.Range("E49").Value = "  +6.30%  "
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "  -6.10%  "
